$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: H11: 17871740 -> 11372887, I11: 17871740 -> 11372887, K11: 17871740 -> 11372887, M11: -17871600 -> -11372747
$ws.Range("H11").Value = 11372887
$ws.Range("I11").Value = 11372887
$ws.Range("K11").Value = 11372887
$ws.Range("M11").Value = -11372747
# Row 46: H46: 189373.86 -> 190073.72, I46: 267039 -> 400508.5, J46: 131125 -> 105899.8, K46: 801117 -> 1201525.5, L46: 393375 -> 317699.4, M46: -800998 -> -1201406.5, N46: -393613 -> -317937.4
$ws.Range("H46").Value = 190073.72
$ws.Range("I46").Value = 400508.5
$ws.Range("J46").Value = 105899.8
$ws.Range("K46").Value = 1201525.5
$ws.Range("L46").Value = 317699.4
$ws.Range("M46").Value = -1201406.5
$ws.Range("N46").Value = -317937.4
# Row 60: H60: 189373.86 -> 190073.72, I60: 267039 -> 400508.5, J60: 131125 -> 105899.8, K60: 801117 -> 1201525.5, L60: 393375 -> 317699.4, M60: -800633 -> -1201041.5, N60: -394343 -> -318667.4
$ws.Range("H60").Value = 190073.72
$ws.Range("I60").Value = 400508.5
$ws.Range("J60").Value = 105899.8
$ws.Range("K60").Value = 1201525.5
$ws.Range("L60").Value = 317699.4
$ws.Range("M60").Value = -1201041.5
$ws.Range("N60").Value = -318667.4
# Row 128: H128: 40000 -> 45390, J128: 0 -> 50780, L128: 0 -> 50780, N128: NEW -> -60740
$ws.Range("H128").Value = 45390
$ws.Range("J128").Value = 50780
$ws.Range("L128").Value = 50780
$ws.Range("N128").Value = -60740
# Row 132: H132: 682247.75 -> 621679.1, I132: 1654.6727 -> 1322.8871, J132: 2884166.5 -> 2884154.8, K132: 4964.0181 -> 3968.6613, L132: 8652499.5 -> 8652464.399999999, M132: -2434.0181 -> -1438.6613, N132: -8657559.5 -> -8657524.399999999
$ws.Range("H132").Value = 621679.1
$ws.Range("I132").Value = 1322.8871
$ws.Range("J132").Value = 2884154.8
$ws.Range("K132").Value = 3968.6613
$ws.Range("L132").Value = 8652464.399999999
$ws.Range("M132").Value = -1438.6613
$ws.Range("N132").Value = -8657524.399999999
# Row 137: H137: 2704446.2 -> 2175371.5, I137: 4546898 -> 3449405, J137: 2183.4 -> 2020.6471, K137: 13640694 -> 10348215, L137: 6550.200000000001 -> 6061.9413, M137: -13638144 -> -10345665, N137: -11650.2 -> -11161.9413
$ws.Range("H137").Value = 2175371.5
$ws.Range("I137").Value = 3449405
$ws.Range("J137").Value = 2020.6471
$ws.Range("K137").Value = 10348215
$ws.Range("L137").Value = 6061.9413
$ws.Range("M137").Value = -10345665
$ws.Range("N137").Value = -11161.9413
# Row 138: H138: 1874636.8 -> 1738003.4, I138: 1266.85 -> 1139.5869, J138: 3403918.2 -> 3335918, K138: 3800.55 -> 3418.7607, L138: 10211754.6 -> 10007754, M138: 1339.45 -> 1721.2393, N138: -10222034.6 -> -10018034
$ws.Range("H138").Value = 1738003.4
$ws.Range("I138").Value = 1139.5869
$ws.Range("J138").Value = 3335918
$ws.Range("K138").Value = 3418.7607
$ws.Range("L138").Value = 10007754
$ws.Range("M138").Value = 1721.2393
$ws.Range("N138").Value = -10018034
# Row 141: H141: 2991.9546 -> 2321.1333, I141: 2420.1428 -> 1988.4445, J141: 15000 -> 5315.3335, K141: 7260.428400000001 -> 5965.333500000001, L141: 45000 -> 15946.0005, M141: -2080.428400000001 -> -785.3335000000006, N141: -55360 -> -26306.0005
$ws.Range("H141").Value = 2321.1333
$ws.Range("I141").Value = 1988.4445
$ws.Range("J141").Value = 5315.3335
$ws.Range("K141").Value = 5965.333500000001
$ws.Range("L141").Value = 15946.0005
$ws.Range("M141").Value = -785.3335000000006
$ws.Range("N141").Value = -26306.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32: 1310.97 -> 864.53, I32: 876.4286 -> 726.19354, J32: 3592.3125 -> 2702.4285, K32: 876.4286 -> 726.19354, L32: 3592.3125 -> 2702.4285, M32: -589.4286 -> -439.19354, N32: -4166.3125 -> -3276.4285
$ws.Range("H32").Value = 864.53
$ws.Range("I32").Value = 726.19354
$ws.Range("J32").Value = 2702.4285
$ws.Range("K32").Value = 726.19354
$ws.Range("L32").Value = 2702.4285
$ws.Range("M32").Value = -439.19354
$ws.Range("N32").Value = -3276.4285
# Row 61: H61: 19270244 -> 14315103, I61: 20855152 -> 16145994, J61: 251353.5 -> 125703.5, K61: 20855152 -> 16145994, L61: 251353.5 -> 125703.5, M61: -20854940 -> -16145782, N61: -251777.5 -> -126127.5
$ws.Range("H61").Value = 14315103
$ws.Range("I61").Value = 16145994
$ws.Range("J61").Value = 125703.5
$ws.Range("K61").Value = 16145994
$ws.Range("L61").Value = 125703.5
$ws.Range("M61").Value = -16145782
$ws.Range("N61").Value = -126127.5
# Row 63: H63: 2011.6666 -> 2031.25, I63: 1975.625 -> 1992.8572, K63: 1975.625 -> 1992.8572, M63: -1289.625 -> -1306.8572
$ws.Range("H63").Value = 2031.25
$ws.Range("I63").Value = 1992.8572
$ws.Range("K63").Value = 1992.8572
$ws.Range("M63").Value = -1306.8572
# Row 66: H66: 2011.6666 -> 2031.25, I66: 1975.625 -> 1992.8572, K66: 9878.125 -> 9964.286, M66: -6446.125 -> -6532.286
$ws.Range("H66").Value = 2031.25
$ws.Range("I66").Value = 1992.8572
$ws.Range("K66").Value = 9964.286
$ws.Range("M66").Value = -6532.286
# Row 97: H97: 1786724.6 -> 2155998, I97: 2404890.5 -> 2976973.5, J97: 912.3333 -> 937.5, K97: 2404890.5 -> 2976973.5, L97: 912.3333 -> 937.5, M97: -2404394.5 -> -2976477.5, N97: -1904.3333 -> -1929.5
$ws.Range("H97").Value = 2155998
$ws.Range("I97").Value = 2976973.5
$ws.Range("J97").Value = 937.5
$ws.Range("K97").Value = 2976973.5
$ws.Range("L97").Value = 937.5
$ws.Range("M97").Value = -2976477.5
$ws.Range("N97").Value = -1929.5
# Row 125: H125: 60000 -> 53551.723, J125: 60000 -> 53551.723, L125: 60000 -> 53551.723, N125: -69840 -> -63391.723
$ws.Range("H125").Value = 53551.723
$ws.Range("J125").Value = 53551.723
$ws.Range("L125").Value = 53551.723
$ws.Range("N125").Value = -63391.723
# Row 132: H132: 35582.566 -> 30095.438, I132: 21578.084 -> 17882.896, J132: 91600.5 -> 84582.16, K132: 64734.25199999999 -> 53648.688, L132: 274801.5 -> 253746.48, M132: -62204.25199999999 -> -51118.688, N132: -279861.5 -> -258806.48
$ws.Range("H132").Value = 30095.438
$ws.Range("I132").Value = 17882.896
$ws.Range("J132").Value = 84582.16
$ws.Range("K132").Value = 53648.688
$ws.Range("L132").Value = 253746.48
$ws.Range("M132").Value = -51118.688
$ws.Range("N132").Value = -258806.48
# Row 136: H136: 19270244 -> 14315103, I136: 20855152 -> 16145994, J136: 251353.5 -> 125703.5, K136: 62565456 -> 48437982, L136: 754060.5 -> 377110.5, M136: -62562906 -> -48435432, N136: -759160.5 -> -382210.5
$ws.Range("H136").Value = 14315103
$ws.Range("I136").Value = 16145994
$ws.Range("J136").Value = 125703.5
$ws.Range("K136").Value = 48437982
$ws.Range("L136").Value = 377110.5
$ws.Range("M136").Value = -48435432
$ws.Range("N136").Value = -382210.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82: H82: 28119.2 -> 23682.666, I82: 9875 -> 7083.3335, K82: 9875 -> 7083.3335, M82: -9492 -> -6700.3335
$ws.Range("H82").Value = 23682.666
$ws.Range("I82").Value = 7083.3335
$ws.Range("K82").Value = 7083.3335
$ws.Range("M82").Value = -6700.3335
# Row 85: H85: 28119.2 -> 23682.666, I85: 9875 -> 7083.3335, K85: 9875 -> 7083.3335, M85: -8549 -> -5757.3335
$ws.Range("H85").Value = 23682.666
$ws.Range("I85").Value = 7083.3335
$ws.Range("K85").Value = 7083.3335
$ws.Range("M85").Value = -5757.3335

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31: 2403.6345 -> 2046.8438, I31: 1154.8857 -> 1006.6667, J31: 4974.5884 -> 4510.421, K31: 1154.8857 -> 1006.6667, L31: 4974.5884 -> 4510.421, M31: -859.8857 -> -711.6667, N31: -5564.5884 -> -5100.421
$ws.Range("H31").Value = 2046.8438
$ws.Range("I31").Value = 1006.6667
$ws.Range("J31").Value = 4510.421
$ws.Range("K31").Value = 1006.6667
$ws.Range("L31").Value = 4510.421
$ws.Range("M31").Value = -711.6667
$ws.Range("N31").Value = -5100.421
# Row 34: H34: 2403.6345 -> 2046.8438, I34: 1154.8857 -> 1006.6667, J34: 4974.5884 -> 4510.421, K34: 1154.8857 -> 1006.6667, L34: 4974.5884 -> 4510.421, M34: -952.8857 -> -804.6667, N34: -5378.5884 -> -4914.421
$ws.Range("H34").Value = 2046.8438
$ws.Range("I34").Value = 1006.6667
$ws.Range("J34").Value = 4510.421
$ws.Range("K34").Value = 1006.6667
$ws.Range("L34").Value = 4510.421
$ws.Range("M34").Value = -804.6667
$ws.Range("N34").Value = -4914.421
# Row 41: H41: 12000 -> 12750
$ws.Range("H41").Value = 12750
# Row 58: H58: 21278130 -> 16394698, I58: 23257178 -> 18182950, J58: 3375.25 -> 2383.5, K58: 23257178 -> 18182950, L58: 3375.25 -> 2383.5, M58: -23256975 -> -18182747, N58: -3781.25 -> -2789.5
$ws.Range("H58").Value = 16394698
$ws.Range("I58").Value = 18182950
$ws.Range("J58").Value = 2383.5
$ws.Range("K58").Value = 18182950
$ws.Range("L58").Value = 2383.5
$ws.Range("M58").Value = -18182747
$ws.Range("N58").Value = -2789.5
# Row 94: H94: 4921.846 -> 4895.846, I94: 23100 -> 15734, J94: 1616.7273 -> 1644.4, K94: 23100 -> 15734, L94: 1616.7273 -> 1644.4, M94: -22649 -> -15283, N94: -2518.7273 -> -2546.4
$ws.Range("H94").Value = 4895.846
$ws.Range("I94").Value = 15734
$ws.Range("J94").Value = 1644.4
$ws.Range("K94").Value = 15734
$ws.Range("L94").Value = 1644.4
$ws.Range("M94").Value = -15283
$ws.Range("N94").Value = -2546.4
# Row 99: H99: 3218.182 -> 3017.8333, I99: 1550 -> 1520, J99: 4171.4287 -> 4087.7144, K99: 1550 -> 1520, L99: 4171.4287 -> 4087.7144, M99: -52 -> -22, N99: -7167.4287 -> -7083.7144
$ws.Range("H99").Value = 3017.8333
$ws.Range("I99").Value = 1520
$ws.Range("J99").Value = 4087.7144
$ws.Range("K99").Value = 1520
$ws.Range("L99").Value = 4087.7144
$ws.Range("M99").Value = -22
$ws.Range("N99").Value = -7083.7144
# Row 126: H126: 3218.182 -> 3017.8333, I126: 1550 -> 1520, J126: 4171.4287 -> 4087.7144, K126: 4650 -> 4560, L126: 12514.2861 -> 12263.1432, M126: -2180 -> -2090, N126: -17454.2861 -> -17203.1432
$ws.Range("H126").Value = 3017.8333
$ws.Range("I126").Value = 1520
$ws.Range("J126").Value = 4087.7144
$ws.Range("K126").Value = 4560
$ws.Range("L126").Value = 12263.1432
$ws.Range("M126").Value = -2090
$ws.Range("N126").Value = -17203.1432
# Row 129: H129: 31598.8 -> 26608.666, J129: 31598.8 -> 26608.666, L129: 31598.8 -> 26608.666, N129: -41598.8 -> -36608.666
$ws.Range("H129").Value = 26608.666
$ws.Range("J129").Value = 26608.666
$ws.Range("L129").Value = 26608.666
$ws.Range("N129").Value = -36608.666
# Row 134: H134: 18233.35 -> 17684.139, I134: 1038.34 -> 1021.27454, J134: 84368 -> 78384.57000000001, K134: 3115.02 -> 3063.82362, L134: 253104 -> 235153.71, M134: -580.0199999999995 -> -528.8236200000001, N134: -258174 -> -240223.71
$ws.Range("H134").Value = 17684.139
$ws.Range("I134").Value = 1021.27454
$ws.Range("J134").Value = 78384.57000000001
$ws.Range("K134").Value = 3063.82362
$ws.Range("L134").Value = 235153.71
$ws.Range("M134").Value = -528.8236200000001
$ws.Range("N134").Value = -240223.71
# Row 136: H136: 21278130 -> 16394698, I136: 23257178 -> 18182950, J136: 3375.25 -> 2383.5, K136: 69771534 -> 54548850, L136: 10125.75 -> 7150.5, M136: -69768984 -> -54546300, N136: -15225.75 -> -12250.5
$ws.Range("H136").Value = 16394698
$ws.Range("I136").Value = 18182950
$ws.Range("J136").Value = 2383.5
$ws.Range("K136").Value = 54548850
$ws.Range("L136").Value = 7150.5
$ws.Range("M136").Value = -54546300
$ws.Range("N136").Value = -12250.5

$ws = $wb.Worksheets.Item("CUL")
# Row 116: H116: 167854 -> 149248, I116: 333300 -> 167450, J116: 112705.336 -> 134686.4, K116: 999900 -> 502350, L116: 338116.008 -> 404059.2, M116: -996458 -> -498908, N116: -345000.008 -> -410943.2
$ws.Range("H116").Value = 149248
$ws.Range("I116").Value = 167450
$ws.Range("J116").Value = 134686.4
$ws.Range("K116").Value = 502350
$ws.Range("L116").Value = 404059.2
$ws.Range("M116").Value = -498908
$ws.Range("N116").Value = -410943.2

$ws = $wb.Worksheets.Item("GSM")
# Row 102: H102: 657.6923 -> 652.7857, I102: 645.8333 -> 681.9, J102: 800 -> 580, K102: 645.8333 -> 681.9, L102: 800 -> 580, M102: 976.1667 -> 940.1, N102: -4044 -> -3824
$ws.Range("H102").Value = 652.7857
$ws.Range("I102").Value = 681.9
$ws.Range("J102").Value = 580
$ws.Range("K102").Value = 681.9
$ws.Range("L102").Value = 580
$ws.Range("M102").Value = 940.1
$ws.Range("N102").Value = -3824

$ws = $wb.Worksheets.Item("LTW")
# Row 132: H132: 93545.17999999999 -> 57383.777, I132: 4200 -> 1722, J132: 144599.58 -> 252200, K132: 12600 -> 5166, L132: 433798.74 -> 756600, M132: -10070 -> -2636, N132: -438858.74 -> -761660
$ws.Range("H132").Value = 57383.777
$ws.Range("I132").Value = 1722
$ws.Range("J132").Value = 252200
$ws.Range("K132").Value = 5166
$ws.Range("L132").Value = 756600
$ws.Range("M132").Value = -2636
$ws.Range("N132").Value = -761660

$ws = $wb.Worksheets.Item("WVR")
# Row 132: H132: 39285.77 -> 48413.26, I132: 22514.674 -> 36383.695, J132: 167864.17 -> 123263.89, K132: 67544.022 -> 109151.085, L132: 503592.51 -> 369791.67, M132: -65014.022 -> -106621.085, N132: -508652.51 -> -374851.67
$ws.Range("H132").Value = 48413.26
$ws.Range("I132").Value = 36383.695
$ws.Range("J132").Value = 123263.89
$ws.Range("K132").Value = 109151.085
$ws.Range("L132").Value = 369791.67
$ws.Range("M132").Value = -106621.085
$ws.Range("N132").Value = -374851.67
